# maj TP06 cours 7
# Fill the "G" column (missing values) for the first set of rows with 1,
# and update the active selection/scroll position on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose G cell is currently empty and must be set to 1.
$rows = @(3,4,5,6,7,8,9,10,11,12,13,16,17)
foreach ($r in $rows) {
    $ws.Range("G$r").Value = 1
}

# Update the selection / view: the sheet no longer scrolls to A6, and the
# active cell / selection becomes G9 instead of F6.
$ws.Range("A1").Select()
$ws.Range("G9").Select()
